$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1561.5883
$ws.Range("J19").Value = 1139.2858
$ws.Range("L19").Value = 1139.2858
$ws.Range("N19").Value = -1489.2858
$ws.Range("H62").Value = 2998
$ws.Range("J62").Value = 2998
$ws.Range("L62").Value = 2998
$ws.Range("N62").Value = -4246
$ws.Range("H64").Value = 4375
$ws.Range("I64").Value = 5000
$ws.Range("K64").Value = 5000
$ws.Range("M64").Value = -4752
$ws.Range("H65").Value = 2998
$ws.Range("J65").Value = 2998
$ws.Range("L65").Value = 14990
$ws.Range("N65").Value = -21230
$ws.Range("H67").Value = 4375
$ws.Range("I67").Value = 5000
$ws.Range("K67").Value = 5000
$ws.Range("M67").Value = -4142
$ws.Range("H100").Value = 1001847.3
$ws.Range("I100").Value = 2002196
$ws.Range("J100").Value = 1498.6
$ws.Range("K100").Value = 2002196
$ws.Range("L100").Value = 1498.6
$ws.Range("M100").Value = -2001655
$ws.Range("N100").Value = -2580.6
$ws.Range("H107").Value = 564.6
$ws.Range("I107").Value = 656
$ws.Range("J107").Value = 199
$ws.Range("K107").Value = 656
$ws.Range("L107").Value = 199
$ws.Range("M107").Value = 1264
$ws.Range("N107").Value = -4039
$ws.Range("H112").Value = 1711.3529
$ws.Range("J112").Value = 1913.7858
$ws.Range("L112").Value = 5741.357400000001
$ws.Range("N112").Value = -7957.357400000001
$ws.Range("H115").Value = 9000
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H116").Value = 6983.143
$ws.Range("J116").Value = 6975.6
$ws.Range("L116").Value = 6975.6
$ws.Range("N116").Value = -13859.6
$ws.Range("H118").Value = 1000
$ws.Range("I118").Value = 1000
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 3000
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -1343
$ws.Range("N118").ClearContents()
$ws.Range("H120").Value = 100000
$ws.Range("J120").Value = 100000
$ws.Range("L120").Value = 100000
$ws.Range("N120").Value = -109676
$ws.Range("H121").Value = 1945.6666
$ws.Range("J121").Value = 1945.6666
$ws.Range("L121").Value = 5836.9998
$ws.Range("N121").Value = -9330.9998
$ws.Range("H132").Value = 1817
$ws.Range("I132").Value = 1146.8948
$ws.Range("K132").Value = 3440.6844
$ws.Range("M132").Value = -910.6844000000001
$ws.Range("H135").Value = 2076.2
$ws.Range("I135").Value = 2076.2
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 18685.8
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -16150.8
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 1742.4
$ws.Range("I137").Value = 1268.5
$ws.Range("J137").Value = 2690.2
$ws.Range("K137").Value = 3805.5
$ws.Range("L137").Value = 8070.599999999999
$ws.Range("M137").Value = -1255.5
$ws.Range("N137").Value = -13170.6
$ws.Range("H138").Value = 2680.3057
$ws.Range("J138").Value = 3052.56
$ws.Range("L138").Value = 9157.68
$ws.Range("N138").Value = -19437.68
$ws.Range("H141").Value = 4331.3335
$ws.Range("I141").Value = 4406.909
$ws.Range("J141").Value = 3500
$ws.Range("K141").Value = 13220.727
$ws.Range("L141").Value = 10500
$ws.Range("M141").Value = -8040.726999999999
$ws.Range("N141").Value = -20860

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1234.3572
$ws.Range("I2").Value = 1327.091
$ws.Range("J2").Value = 894.3333
$ws.Range("K2").Value = 1327.091
$ws.Range("L2").Value = 894.3333
$ws.Range("M2").Value = -1214.091
$ws.Range("N2").Value = -1120.3333
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 28562530
$ws.Range("I74").Value = 39985740
$ws.Range("K74").Value = 39985740
$ws.Range("M74").Value = -39984866
$ws.Range("H77").Value = 28562530
$ws.Range("I77").Value = 39985740
$ws.Range("K77").Value = 199928700
$ws.Range("M77").Value = -199924332
$ws.Range("H116").Value = 1234.3572
$ws.Range("I116").Value = 1327.091
$ws.Range("J116").Value = 894.3333
$ws.Range("K116").Value = 1327.091
$ws.Range("L116").Value = 894.3333
$ws.Range("M116").Value = 966.9090000000001
$ws.Range("N116").Value = -5482.3333
$ws.Range("H122").Value = 4578
$ws.Range("I122").Value = 4690.7144
$ws.Range("K122").Value = 14072.1432
$ws.Range("M122").Value = -11622.1432
$ws.Range("H132").Value = 2771.3215
$ws.Range("I132").Value = 2115.842
$ws.Range("K132").Value = 6347.526
$ws.Range("M132").Value = -3817.526
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1234.3572
$ws.Range("I3").Value = 1327.091
$ws.Range("J3").Value = 894.3333
$ws.Range("K3").Value = 1327.091
$ws.Range("L3").Value = 894.3333
$ws.Range("M3").Value = -1213.091
$ws.Range("N3").Value = -1122.3333
$ws.Range("H94").Value = 506.58334
$ws.Range("I94").Value = 370.81818
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 370.81818
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = 80.18182000000002
$ws.Range("N94").Value = -2902
$ws.Range("H134").Value = 2107.862
$ws.Range("I134").Value = 1901.4783
$ws.Range("K134").Value = 5704.4349
$ws.Range("M134").Value = -3169.4349

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").ClearContents()
$ws.Range("H41").Value = 39618.332
$ws.Range("J41").Value = 39927.5
$ws.Range("L41").Value = 39927.5
$ws.Range("N41").Value = -40783.5
$ws.Range("H62").Value = 4900
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 4900
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H105").Value = 1705
$ws.Range("I105").Value = 1705
$ws.Range("K105").Value = 1705
$ws.Range("M105").Value = 42
$ws.Range("H132").Value = 2934.55
$ws.Range("I132").Value = 2180.6924
$ws.Range("K132").Value = 6542.0772
$ws.Range("M132").Value = -4012.0772
$ws.Range("H134").Value = 4311.375
$ws.Range("I134").Value = 4332
$ws.Range("J134").Value = 4249.5
$ws.Range("K134").Value = 12996
$ws.Range("L134").Value = 12748.5
$ws.Range("M134").Value = -10461
$ws.Range("N134").Value = -17818.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 188778780
$ws.Range("I4").Value = 104079180
$ws.Range("K4").Value = 312237540
$ws.Range("M4").Value = -312237428
$ws.Range("H23").Value = 267.5
$ws.Range("J23").Value = 383.33334
$ws.Range("L23").Value = 1150.00002
$ws.Range("N23").Value = -1620.00002
$ws.Range("H131").Value = 1499.5
$ws.Range("I131").Value = 1000
$ws.Range("J131").Value = 1999
$ws.Range("K131").Value = 3000
$ws.Range("L131").Value = 5997
$ws.Range("M131").Value = 2040
$ws.Range("N131").Value = -16077

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1466.1538
$ws.Range("I102").Value = 1345
$ws.Range("K102").Value = 1345
$ws.Range("M102").Value = 277
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4864

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4302.6665
$ws.Range("I132").Value = 4496.5
$ws.Range("K132").Value = 13489.5
$ws.Range("M132").Value = -10959.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3574.1428
$ws.Range("I122").Value = 3304
$ws.Range("K122").Value = 9912
$ws.Range("M122").Value = -7462
$ws.Range("H132").Value = 3609.3333
$ws.Range("I132").Value = 3492.0833
$ws.Range("J132").Value = 3726.5833
$ws.Range("K132").Value = 10476.2499
$ws.Range("L132").Value = 11179.7499
$ws.Range("M132").Value = -7946.249899999999
$ws.Range("N132").Value = -16239.7499
